# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF. Copy the formatting from an
# existing header cell (bold font + thin border, style index 1) so the new
# header cells match the look of the rest of the header row, then set text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2 through 54 all get the same season record values.
$ws.Range("AD2:AD54").Value = 77
$ws.Range("AE2:AE54").Value = 85
$ws.Range("AF2:AF54").Value = 0
